$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows (162, 163) by copying the formatting of the last
# existing data row (161) so the date/number styles match the rest of the
# table, then overwrite the copied values with the real data for
# 2025-11-20 (Excel serial date 45981).
$ws.Rows.Item(161).Copy()
$ws.Rows.Item(162).Insert()
$ws.Rows.Item(161).Copy()
$ws.Rows.Item(163).Insert()

# Row 162: 四方坪站充电量(kw)
$ws.Cells.Item(162, 1).Value = 45981
$ws.Cells.Item(162, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(162, 3).Value = 580.41199999999992
$ws.Cells.Item(162, 4).Value = 1073.7449999999999
$ws.Cells.Item(162, 5).Value = 557.44999999999993
$ws.Cells.Item(162, 6).Value = 433.66300000000001
$ws.Cells.Item(162, 7).Value = 300.21100000000001
$ws.Cells.Item(162, 8).Value = 575.346
$ws.Cells.Item(162, 9).Value = 470.05599999999993
$ws.Cells.Item(162, 10).Value = 182.21899999999999
$ws.Cells.Item(162, 11).Value = 122.34
$ws.Cells.Item(162, 12).Value = 78.63
$ws.Cells.Item(162, 13).Value = 278.8
$ws.Cells.Item(162, 14).Value = 257.36200000000002
$ws.Cells.Item(162, 15).Value = 721.62500000000011
$ws.Cells.Item(162, 16).Value = 1451.1109999999999
$ws.Cells.Item(162, 17).Value = 565.21300000000008
$ws.Cells.Item(162, 18).Value = 433.66800000000006
$ws.Cells.Item(162, 19).Value = 326.17900000000003
$ws.Cells.Item(162, 20).Value = 179.739
$ws.Cells.Item(162, 21).Value = 154.05399999999997
$ws.Cells.Item(162, 22).Value = 199.684
$ws.Cells.Item(162, 23).Value = 155.97999999999999
$ws.Cells.Item(162, 24).Value = 87.42
$ws.Cells.Item(162, 25).Value = 106.92
$ws.Cells.Item(162, 26).Value = 49.26

# Row 163: 高岭站充电量(kw)
$ws.Cells.Item(163, 1).Value = 45981
$ws.Cells.Item(163, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(163, 3).Value = 396.28499999999997
$ws.Cells.Item(163, 4).Value = 364.96800000000002
$ws.Cells.Item(163, 5).Value = 77.546999999999997
$ws.Cells.Item(163, 6).Value = 74.992999999999995
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 141.42099999999999
$ws.Cells.Item(163, 9).Value = 140
$ws.Cells.Item(163, 10).Value = 120.104
$ws.Cells.Item(163, 11).Value = 292.18799999999999
$ws.Cells.Item(163, 12).Value = 110.786
$ws.Cells.Item(163, 13).Value = 208.21399999999997
$ws.Cells.Item(163, 14).Value = 159.197
$ws.Cells.Item(163, 15).Value = 434.38300000000004
$ws.Cells.Item(163, 16).Value = 419.37799999999993
$ws.Cells.Item(163, 17).Value = 146.82499999999999
$ws.Cells.Item(163, 18).Value = 217.79900000000001
$ws.Cells.Item(163, 19).Value = 214.24200000000002
$ws.Cells.Item(163, 20).Value = 45.936999999999998
$ws.Cells.Item(163, 21).Value = 120.33600000000001
$ws.Cells.Item(163, 22).Value = 71.989000000000004
$ws.Cells.Item(163, 23).Value = 0
$ws.Cells.Item(163, 24).Value = 59.983999999999995
$ws.Cells.Item(163, 25).Value = 179.839
$ws.Cells.Item(163, 26).Value = 87.644999999999996
